# Apply commit "100% test coverage for researcher, including all memoization
# buffers" — adds three new timing sheets (copies of the previous
# "Created update_history buffer" sheet, each with its own input numbers and
# a "Note:" annotation), and tidies up the selection/tab state of the sheet
# they were copied from.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Fixed playerInfo appends" - straight copy of the previous sheet with
#    new timing numbers, no Note column yet.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("Created update_history buffer")
$src.Copy($null, $src)
$s11 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s11.Name = "Fixed playerInfo appends"

$s11.Range("B2").Value = 7.8789999999999996
$s11.Range("C2").Value = 7.7649999999999997
$s11.Range("D2").Value = 8.0229999999999997
$s11.Range("G2").Value = 41809

$s11.Range("B3").Value = 7.1550000000000002
$s11.Range("C3").Value = 7.0510000000000002
$s11.Range("D3").Value = 7.2960000000000003

$s11.Range("B4").Value = 0.72
$s11.Range("C4").Value = 0.70799999999999996
$s11.Range("D4").Value = 0.72499999999999998

$s11.Range("A1:G6").Select()
$s11.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 2. "Switched to iTerm2" - copy of sheet 11, new numbers, plus the first
#    "Note:" / explanation pair in columns H/I.
# ---------------------------------------------------------------------
$s11.Copy($null, $s11)
$s12 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s12.Name = "Switched to iTerm2"

$s12.Range("B2").Value = 8.0250000000000004
$s12.Range("C2").Value = 7.8570000000000002
$s12.Range("D2").Value = 7.8369999999999997
$s12.Range("G2").Value = 41809

$s12.Range("B3").Value = 7.1479999999999997
$s12.Range("C3").Value = 7.1219999999999999
$s12.Range("D3").Value = 7.1239999999999997

$s12.Range("B4").Value = 0.72799999999999998
$s12.Range("C4").Value = 0.73
$s12.Range("D4").Value = 0.70799999999999996

$s12.Range("H1").Value = "Note:"
$s12.Range("H1").Font.Bold = $true
$s12.Range("I1").Font.Bold = $true
$s12.Range("H2").Value = "Switched from testing in terminal"
$s12.Range("H3").Value = "to testing in iTerm2"

$s12.Range("H1").Select()
$s12.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. "Change to boxscoreBuffer" - copy of sheet 12, new numbers, and its
#    own Note about the boxscore-buffer dictionary lookups. This ends up
#    the last / active sheet.
# ---------------------------------------------------------------------
$s12.Copy($null, $s12)
$s13 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s13.Name = "Change to boxscoreBuffer"

$s13.Range("B2").Value = 7.9889999999999999
$s13.Range("C2").Value = 7.7809999999999997
$s13.Range("D2").Value = 7.9649999999999999
$s13.Range("G2").Value = 41809

$s13.Range("B3").Value = 7.2560000000000002
$s13.Range("C3").Value = 7.0609999999999999
$s13.Range("D3").Value = 7.2359999999999998

$s13.Range("B4").Value = 0.72799999999999998
$s13.Range("C4").Value = 0.71599999999999997
$s13.Range("D4").Value = 0.72399999999999998

# This sheet keeps only the Note: label in H1 (no I1) and a single note in
# H2 - remove the leftover "to testing in iTerm2" cell carried over by Copy.
$s13.Range("I1").Clear()
$s13.Range("H2").Value = "reduced dictionary lookups in boxscore buffer by 33%ish. "
$s13.Range("H3").Clear()

$s13.Range("H2").Select()
$s13.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4. Tidy up the sheet the new ones were copied from: it is no longer the
#    selected tab, and its lingering D5 selection is replaced with the
#    whole-table selection the other sheets use.
# ---------------------------------------------------------------------
$src.Range("A1:G6").Select()
$src.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5. Make the newest sheet ("Change to boxscoreBuffer") the active tab.
# ---------------------------------------------------------------------
$s13.Activate()
